{"js": "// Apply the 2025-12-21 -> 2025-12-22 date/weekday update and refresh the\n// 3-digit x 1-digit multiplication problems in the practice table.\n// Each (oldText, newText) pair is a unique run of text in the document, so\n// a direct search + replace-in-place keeps every run's formatting intact.\nconst replacements = [\n  [\"2025-12-21 Sunday\", \"2025-12-22 Monday\"],\n  [\"759\u00d79=\", \"254\u00d75=\"],\n  [\"827\u00d79=\", \"921\u00d72=\"],\n  [\"505\u00d79=\", \"933\u00d74=\"],\n  [\"774\u00d75=\", \"435\u00d79=\"],\n  [\"120\u00d78=\", \"965\u00d73=\"],\n  [\"692\u00d79=\", \"963\u00d78=\"],\n  [\"306\u00d74=\", \"746\u00d72=\"],\n  [\"611\u00d74=\", \"917\u00d75=\"],\n  [\"240\u00d74=\", \"289\u00d78=\"],\n  [\"870\u00d77=\", \"584\u00d78=\"],\n  [\"323\u00d74=\", \"619\u00d77=\"],\n  [\"596\u00d76=\", \"982\u00d75=\"],\n  [\"656\u00d76=\", \"551\u00d74=\"],\n  [\"529\u00d79=\", \"627\u00d79=\"],\n  [\"947\u00d73=\", \"951\u00d74=\"],\n  [\"845\u00d79=\", \"441\u00d78=\"],\n  [\"331\u00d73=\", \"436\u00d78=\"],\n  [\"354\u00d73=\", \"511\u00d73=\"],\n  [\"975\u00d72=\", \"699\u00d72=\"],\n  [\"104\u00d72=\", \"992\u00d78=\"],\n  [\"793\u00d78=\", \"441\u00d77=\"],\n  [\"510\u00d75=\", \"563\u00d74=\"],\n  [\"778\u00d75=\", \"291\u00d78=\"],\n  [\"423\u00d77=\", \"365\u00d75=\"],\n  [\"638\u00d76=\", \"133\u00d72=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the 2025-12-21 -> 2025-12-22 date/weekday update and refresh\n# the 3-digit x 1-digit multiplication problems in the practice table.\n# Each (old, new) pair is a unique piece of text in the document, so a\n# direct Find/Replace keeps every run's formatting (font, size) intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-21 Sunday\", \"2025-12-22 Monday\"),\n    @(\"759\u00d79=\", \"254\u00d75=\"),\n    @(\"827\u00d79=\", \"921\u00d72=\"),\n    @(\"505\u00d79=\", \"933\u00d74=\"),\n    @(\"774\u00d75=\", \"435\u00d79=\"),\n    @(\"120\u00d78=\", \"965\u00d73=\"),\n    @(\"692\u00d79=\", \"963\u00d78=\"),\n    @(\"306\u00d74=\", \"746\u00d72=\"),\n    @(\"611\u00d74=\", \"917\u00d75=\"),\n    @(\"240\u00d74=\", \"289\u00d78=\"),\n    @(\"870\u00d77=\", \"584\u00d78=\"),\n    @(\"323\u00d74=\", \"619\u00d77=\"),\n    @(\"596\u00d76=\", \"982\u00d75=\"),\n    @(\"656\u00d76=\", \"551\u00d74=\"),\n    @(\"529\u00d79=\", \"627\u00d79=\"),\n    @(\"947\u00d73=\", \"951\u00d74=\"),\n    @(\"845\u00d79=\", \"441\u00d78=\"),\n    @(\"331\u00d73=\", \"436\u00d78=\"),\n    @(\"354\u00d73=\", \"511\u00d73=\"),\n    @(\"975\u00d72=\", \"699\u00d72=\"),\n    @(\"104\u00d72=\", \"992\u00d78=\"),\n    @(\"793\u00d78=\", \"441\u00d77=\"),\n    @(\"510\u00d75=\", \"563\u00d74=\"),\n    @(\"778\u00d75=\", \"291\u00d78=\"),\n    @(\"423\u00d77=\", \"365\u00d75=\"),\n    @(\"638\u00d76=\", \"133\u00d72=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
